# Apply the cryptos list update (prices / 1h volume % / a row re-order at 44-45)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $rng = $Sheet.Range($CellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "43.887.08"
Set-TextValue $ws "E2" "  +0.00%  "
Set-TextValue $ws "D3" "2.360.76"
Set-TextValue $ws "E3" "  +0.47%  "
Set-TextValue $ws "E4" "  -0.11%  "
Set-TextValue $ws "E5" "  -0.53%  "
Set-TextValue $ws "D6" "240.23"
Set-TextValue $ws "E6" "  +0.65%  "
Set-TextValue $ws "D7" "74.02"
Set-TextValue $ws "E7" "  +0.81%  "
Set-TextValue $ws "E8" "  -0.05%  "
Set-TextValue $ws "E9" "  +3.12%  "
Set-TextValue $ws "E10" "  +2.36%  "
Set-TextValue $ws "D11" "60.81"
Set-TextValue $ws "E11" "  +6.20%  "
Set-TextValue $ws "D12" "37.03"
Set-TextValue $ws "E12" "  +15.39%  "
Set-TextValue $ws "E13" "  +0.71%  "
Set-TextValue $ws "D14" "7.29"
Set-TextValue $ws "E14" "  +0.31%  "
Set-TextValue $ws "D15" "16.29"
Set-TextValue $ws "E15" "  -1.33%  "
Set-TextValue $ws "D16" "0.919"
Set-TextValue $ws "E16" "  +2.76%  "
Set-TextValue $ws "D17" "2.362.24"
Set-TextValue $ws "E17" "  -0.06%  "
Set-TextValue $ws "D18" "43.827.09"
Set-TextValue $ws "E18" "  +0.02%  "
Set-TextValue $ws "E19" "  +1.72%  "
Set-TextValue $ws "D20" "77.96"
Set-TextValue $ws "E20" "  +1.76%  "
Set-TextValue $ws "E21" "  -1.80%  "
Set-TextValue $ws "D22" "251.69"
Set-TextValue $ws "E22" "  -2.11%  "
Set-TextValue $ws "D23" "1.00"
Set-TextValue $ws "E23" "  -0.04%  "
Set-TextValue $ws "D24" "3.78"
Set-TextValue $ws "E24" "  +3.47%  "
Set-TextValue $ws "D25" "1.88"
Set-TextValue $ws "E25" "  -2.56%  "
Set-TextValue $ws "E26" "  +0.61%  "
Set-TextValue $ws "D27" "10.52"
Set-TextValue $ws "E27" "  -1.66%  "
Set-TextValue $ws "E28" "  +0.57%  "
Set-TextValue $ws "D29" "22.35"
Set-TextValue $ws "E29" "  -0.95%  "
Set-TextValue $ws "D30" "175.64"
Set-TextValue $ws "E30" "  +0.06%  "
Set-TextValue $ws "E31" "  +1.18%  "
Set-TextValue $ws "E32" "  -1.41%  "
Set-TextValue $ws "E33" "  -1.44%  "
Set-TextValue $ws "D34" "5.09"
Set-TextValue $ws "E34" "  -1.99%  "
Set-TextValue $ws "D35" "5.39"
Set-TextValue $ws "E35" "  -0.49%  "
Set-TextValue $ws "E36" "  +1.88%  "
Set-TextValue $ws "E37" "  +5.51%  "
Set-TextValue $ws "D38" "2.41"
Set-TextValue $ws "E38" "  +3.12%  "
Set-TextValue $ws "E39" "  +0.64%  "
Set-TextValue $ws "D40" "5.43"
Set-TextValue $ws "E40" "  +14.27%  "
Set-TextValue $ws "D41" "20.36"
Set-TextValue $ws "E41" "  +7.78%  "
Set-TextValue $ws "D42" "65.24"
Set-TextValue $ws "E42" "  +12.83%  "
Set-TextValue $ws "E43" "  -0.50%  "
Set-TextValue $ws "B44" "FraxShare"
Set-TextValue $ws "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws "D44" "9.06"
Set-TextValue $ws "E44" "  +0.69%  "
Set-TextValue $ws "B45" "Cronos"
Set-TextValue $ws "C45" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D45" "0.107"
Set-TextValue $ws "E45" "  -2.64%  "
Set-TextValue $ws "D46" "2.54"
Set-TextValue $ws "E46" "  +1.36%  "
Set-TextValue $ws "E47" "  -0.24%  "
Set-TextValue $ws "E48" "  +0.32%  "
Set-TextValue $ws "E49" "  -0.67%  "
Set-TextValue $ws "D50" "98.29"
Set-TextValue $ws "E50" "  -1.31%  "
Set-TextValue $ws "E51" "  +16.24%  "
